# "Update countries & provincias Spain"
#
# This workbook is a COVID-19 ranking table (countries sorted descending by
# total cases). The commit refreshes the underlying data pull, which causes:
#   1. Three countries to leapfrog their immediate neighbour in the ranking
#      (Argentina over Catar, Nigeria over Rumania, Libia over Nueva
#      Zelanda/Yemen) - i.e. the country *names* in column A for a few rows
#      swap places while the row positions stay fixed.
#   2. Fresh case/recovered/death numbers for several rows (including the
#      rows whose country name just changed, plus a handful of unrelated
#      rows further down the ranking).
#   3. The "last updated" footer timestamp text bumps from 00:02 to 01:19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Country-name swaps caused by the re-ranking ---------------------
# Canada(23) / Catar(24) / Argentina(25) / Egipto(26) -> Argentina jumps in
# front of Catar:
$ws.Cells.Item(24, 1).Value = "Argentina"
$ws.Cells.Item(25, 1).Value = "Catar"

# Barein(49) / Rumania(50) / Nigeria(51) / Suiza(52) -> Nigeria jumps in
# front of Rumania:
$ws.Cells.Item(50, 1).Value = "Nigeria"
$ws.Cells.Item(51, 1).Value = "Rumania"

# Hong Kong(126) / Nueva Zelanda(127) / Yemen(128) / Libia(129) / Suazilandia(130)
# -> Libia jumps in front of Nueva Zelanda (and Yemen):
$ws.Cells.Item(127, 1).Value = "Libia"
$ws.Cells.Item(128, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(129, 1).Value = "Yemen"

# --- 2. Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 01:19"

# --- 3. Refreshed numeric data (Casos totales, Nuevos casos, Casos activos,
#        Recuperados, Casos criticos, Muertes hoy, Muertes) -----------------

# Estados Unidos
$ws.Cells.Item(4, 2).Value = 3542790
$ws.Cells.Item(4, 3).Value = 63307
$ws.Cells.Item(4, 4).Value = 1589660
$ws.Cells.Item(4, 5).Value = 1813989
$ws.Cells.Item(4, 7).Value = 894
$ws.Cells.Item(4, 8).Value = 139141

# Brasil
$ws.Cells.Item(5, 2).Value = 1931204
$ws.Cells.Item(5, 3).Value = 43245
$ws.Cells.Item(5, 5).Value = 643430
$ws.Cells.Item(5, 7).Value = 1341
$ws.Cells.Item(5, 8).Value = 74262

# Colombia
$ws.Cells.Item(22, 2).Value = 159898
$ws.Cells.Item(22, 3).Value = 5621
$ws.Cells.Item(22, 4).Value = 68806
$ws.Cells.Item(22, 5).Value = 85467
$ws.Cells.Item(22, 7).Value = 170
$ws.Cells.Item(22, 8).Value = 5625

# Canada
$ws.Cells.Item(23, 2).Value = 108486
$ws.Cells.Item(23, 3).Value = 331
$ws.Cells.Item(23, 4).Value = 72170
$ws.Cells.Item(23, 5).Value = 27518
$ws.Cells.Item(23, 7).Value = 8
$ws.Cells.Item(23, 8).Value = 8798

# Argentina (now at row 24)
$ws.Cells.Item(24, 2).Value = 106910
$ws.Cells.Item(24, 3).Value = 3645
$ws.Cells.Item(24, 4).Value = 45467
$ws.Cells.Item(24, 5).Value = 59475
$ws.Cells.Item(24, 7).Value = 65
$ws.Cells.Item(24, 8).Value = 1968

# Catar (now at row 25)
$ws.Cells.Item(25, 2).Value = 104533
$ws.Cells.Item(25, 3).Value = 517
$ws.Cells.Item(25, 4).Value = 101160
$ws.Cells.Item(25, 5).Value = 3223
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 150

# Nigeria (now at row 50)
$ws.Cells.Item(50, 2).Value = 33616
$ws.Cells.Item(50, 3).Value = 463
$ws.Cells.Item(50, 4).Value = 13792
$ws.Cells.Item(50, 5).Value = 19070
$ws.Cells.Item(50, 7).Value = 10
$ws.Cells.Item(50, 8).Value = 754

# Rumania (now at row 51)
$ws.Cells.Item(51, 2).Value = 33585
$ws.Cells.Item(51, 3).Value = 637
$ws.Cells.Item(51, 4).Value = 21803
$ws.Cells.Item(51, 5).Value = 9851
$ws.Cells.Item(51, 7).Value = 30
$ws.Cells.Item(51, 8).Value = 1931

# Guatemala
$ws.Cells.Item(54, 2).Value = 30872
$ws.Cells.Item(54, 3).Value = 1130
$ws.Cells.Item(54, 4).Value = 4453
$ws.Cells.Item(54, 5).Value = 25117
$ws.Cells.Item(54, 7).Value = 58
$ws.Cells.Item(54, 8).Value = 1302

# Japon
$ws.Cells.Item(59, 2).Value = 22220
$ws.Cells.Item(59, 3).Value = 352
$ws.Cells.Item(59, 4).Value = 18282
$ws.Cells.Item(59, 5).Value = 2956

# Chequia
$ws.Cells.Item(69, 2).Value = 13341
$ws.Cells.Item(69, 3).Value = 103
$ws.Cells.Item(69, 5).Value = 4545

# Venezuela
$ws.Cells.Item(77, 2).Value = 10010
$ws.Cells.Item(77, 3).Value = 303
$ws.Cells.Item(77, 5).Value = 7243
$ws.Cells.Item(77, 7).Value = 3
$ws.Cells.Item(77, 8).Value = 96

# Noruega
$ws.Cells.Item(78, 2).Value = 9001
$ws.Cells.Item(78, 3).Value = 17
$ws.Cells.Item(78, 5).Value = 610

# Libia (now at row 127)
$ws.Cells.Item(127, 2).Value = 1563
$ws.Cells.Item(127, 3).Value = 51
$ws.Cells.Item(127, 4).Value = 370
$ws.Cells.Item(127, 5).Value = 1151
$ws.Cells.Item(127, 7).Value = 2
$ws.Cells.Item(127, 8).Value = 42

# Nueva Zelanda (now at row 128)
$ws.Cells.Item(128, 2).Value = 1545
$ws.Cells.Item(128, 3).Value = 1
$ws.Cells.Item(128, 4).Value = 1498
$ws.Cells.Item(128, 5).Value = 25
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 22

# Yemen (now at row 129)
$ws.Cells.Item(129, 2).Value = 1516
$ws.Cells.Item(129, 3).Value = 18
$ws.Cells.Item(129, 4).Value = 685
$ws.Cells.Item(129, 5).Value = 402
$ws.Cells.Item(129, 7).Value = 5
$ws.Cells.Item(129, 8).Value = 429
